$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")
$ws.Activate()

# Row 2: replace the old lone value with the new set of figures
$ws.Range("A2").Value = 20000
$ws.Range("B2").Value = 560
$ws.Range("C2").Value = 9000
$ws.Range("D2").Value = 25660

# Row 3 (old 14000/22500/86500) is no longer part of the sheet
$ws.Range("B3:D3").ClearContents()

# Leave the selection where the author left it when re-uploading the file
$ws.Range("G6").Select()
